$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (column G) value, derived from the diff
$gValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 0
    6 = 1
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 2
    24 = 3
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 2
    31 = 3
    32 = 3
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 2
    39 = 2
    40 = 0
    41 = 2
    42 = 3
    43 = 1
    44 = 2
    45 = 1
    46 = 1
    47 = 2
    48 = 2
    49 = 3
    50 = 2
    51 = 2
    52 = 4
    53 = 3
    54 = 1
    55 = 3
    56 = 2
    57 = 1
    58 = 3
    59 = 0
    60 = 1
    61 = 2
    62 = 2
    63 = 1
    65 = 2
    66 = 3
    67 = 2
    68 = 3
    70 = 3
    71 = 2
    72 = 2
    73 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

